$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.870.77"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.211.22"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.42"
$ws.Range("E5").Value = "  +5.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.55"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.97"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.543.37"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.49"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.208.46"
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.858.19"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.32"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.50"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.20"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  -7.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.14"
$ws.Range("E26").Value = "  +7.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.79"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -4.81%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0874"
$ws.Range("E33").Value = "  +9.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.23"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E36").Value = "  +7.24%  "
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.91"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +17.72%  "
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.199"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.31"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.02"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.76"
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.33"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.459"
$ws.Range("E48").Value = "  -4.83%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.431.06"
$ws.Range("E51").Value = "  -1.29%  "
